$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 51.59157666666666
$ws.Range("H2").Value = 154.77473
$ws.Range("I2").Value = 0.2641250550177587
$ws.Range("J2").Value = 0.2641250550177588
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.287928333333333
$ws.Range("N2").Value = 15.863785
$ws.Range("O2").Value = 0.08596514992568145
$ws.Range("P2").Value = 0.08596514992568147
$ws.Range("Q2").Value = 272.8125600170055
$ws.Range("R2").Value = 2455.31304015305
$ws.Range("S2").Value = 0.02270554995373049
$ws.Range("T2").Value = 0.0227055499537305
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 51.59157666666666
$ws.Range("H3").Value = 154.77473
$ws.Range("I3").Value = 0.2641250550177587
$ws.Range("J3").Value = 0.2641250550177588
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 39.18893833333333
$ws.Range("N3").Value = 117.566815
$ws.Range("O3").Value = 0.637089375439711
$ws.Range("P3").Value = 0.637089375439711
$ws.Range("Q3").Value = 2021.819116509439
$ws.Range("R3").Value = 18196.37204858494
$ws.Range("S3").Value = 0.1682712663392432
$ws.Range("T3").Value = 0.1682712663392432
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 51.59157666666666
$ws.Range("H4").Value = 154.77473
$ws.Range("I4").Value = 0.2641250550177587
$ws.Range("J4").Value = 0.2641250550177588
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.067472
$ws.Range("N4").Value = 18.202416
$ws.Range("O4").Value = 0.09863808797519777
$ws.Range("P4").Value = 0.09863808797519778
$ws.Range("Q4").Value = 313.0304468608533
$ws.Range("R4").Value = 2817.274021747679
$ws.Range("S4").Value = 0.02605279041329564
$ws.Range("T4").Value = 0.02605279041329565
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 51.59157666666666
$ws.Range("H5").Value = 154.77473
$ws.Range("I5").Value = 0.2641250550177587
$ws.Range("J5").Value = 0.2641250550177588
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 10.968127
$ws.Range("N5").Value = 32.904381
$ws.Range("O5").Value = 0.1783073866594097
$ws.Range("P5").Value = 0.1783073866594098
$ws.Range("Q5").Value = 565.8629650102366
$ws.Range("R5").Value = 5092.76668509213
$ws.Range("S5").Value = 0.04709544831148937
$ws.Range("T5").Value = 0.04709544831148939
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09891538535728452
$ws.Range("J6").Value = 0.09891538535728453
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.287928333333333
$ws.Range("N6").Value = 15.863785
$ws.Range("O6").Value = 0.08596514992568145
$ws.Range("P6").Value = 0.08596514992568147
$ws.Range("Q6").Value = 102.1688741440111
$ws.Range("R6").Value = 919.5198672961002
$ws.Range("S6").Value = 0.00850327593219552
$ws.Range("T6").Value = 0.008503275932195522
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09891538535728452
$ws.Range("J7").Value = 0.09891538535728453
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 39.18893833333333
$ws.Range("N7").Value = 117.566815
$ws.Range("O7").Value = 0.637089375439711
$ws.Range("P7").Value = 0.637089375439711
$ws.Range("Q7").Value = 757.175486508878
$ws.Range("R7").Value = 6814.579378579901
$ws.Range("S7").Value = 0.06301794107865073
$ws.Range("T7").Value = 0.06301794107865073
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09891538535728452
$ws.Range("J8").Value = 0.09891538535728453
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.067472
$ws.Range("N8").Value = 18.202416
$ws.Range("O8").Value = 0.09863808797519777
$ws.Range("P8").Value = 0.09863808797519778
$ws.Range("Q8").Value = 117.2305568577067
$ws.Range("R8").Value = 1055.07501171936
$ws.Range("S8").Value = 0.00975682448297242
$ws.Range("T8").Value = 0.009756824482972421
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09891538535728452
$ws.Range("J9").Value = 0.09891538535728453
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 10.968127
$ws.Range("N9").Value = 32.904381
$ws.Range("O9").Value = 0.1783073866594097
$ws.Range("P9").Value = 0.1783073866594098
$ws.Range("Q9").Value = 211.9168635464734
$ws.Range("R9").Value = 1907.25177191826
$ws.Range("S9").Value = 0.01763734386346585
$ws.Range("T9").Value = 0.01763734386346585
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 112.3724673333333
$ws.Range("H10").Value = 337.117402
$ws.Range("I10").Value = 0.5752951554216499
$ws.Range("J10").Value = 0.57529515542165
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.287928333333333
$ws.Range("N10").Value = 15.863785
$ws.Range("O10").Value = 0.08596514992568145
$ws.Range("P10").Value = 0.08596514992568147
$ws.Range("Q10").Value = 594.2175538985077
$ws.Range("R10").Value = 5347.957985086569
$ws.Range("S10").Value = 0.04945533428734034
$ws.Range("T10").Value = 0.04945533428734036
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 112.3724673333333
$ws.Range("H11").Value = 337.117402
$ws.Range("I11").Value = 0.5752951554216499
$ws.Range("J11").Value = 0.57529515542165
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 39.18893833333333
$ws.Range("N11").Value = 117.566815
$ws.Range("O11").Value = 0.637089375439711
$ws.Range("P11").Value = 0.637089375439711
$ws.Range("Q11").Value = 4403.757692690514
$ws.Range("R11").Value = 39633.81923421462
$ws.Range("S11").Value = 0.3665144312610704
$ws.Range("T11").Value = 0.3665144312610705
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 112.3724673333333
$ws.Range("H12").Value = 337.117402
$ws.Range("I12").Value = 0.5752951554216499
$ws.Range("J12").Value = 0.57529515542165
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.067472
$ws.Range("N12").Value = 18.202416
$ws.Range("O12").Value = 0.09863808797519777
$ws.Range("P12").Value = 0.09863808797519778
$ws.Range("Q12").Value = 681.8167991159145
$ws.Range("R12").Value = 6136.351192043231
$ws.Range("S12").Value = 0.05674601415218578
$ws.Range("T12").Value = 0.0567460141521858
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 112.3724673333333
$ws.Range("H13").Value = 337.117402
$ws.Range("I13").Value = 0.5752951554216499
$ws.Range("J13").Value = 0.57529515542165
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 10.968127
$ws.Range("N13").Value = 32.904381
$ws.Range("O13").Value = 0.1783073866594097
$ws.Range("P13").Value = 0.1783073866594098
$ws.Range("Q13").Value = 1232.515493015351
$ws.Range("R13").Value = 11092.63943713816
$ws.Range("S13").Value = 0.1025793757210533
$ws.Range("T13").Value = 0.1025793757210534
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 12.044915
$ws.Range("H14").Value = 36.134745
$ws.Range("I14").Value = 0.06166440420330686
$ws.Range("J14").Value = 0.06166440420330688
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.287928333333333
$ws.Range("N14").Value = 15.863785
$ws.Range("O14").Value = 0.08596514992568145
$ws.Range("P14").Value = 0.08596514992568147
$ws.Range("Q14").Value = 63.69264730109168
$ws.Range("R14").Value = 573.2338257098251
$ws.Range("S14").Value = 0.005300989752415096
$ws.Range("T14").Value = 0.005300989752415098
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 12.044915
$ws.Range("H15").Value = 36.134745
$ws.Range("I15").Value = 0.06166440420330686
$ws.Range("J15").Value = 0.06166440420330688
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 39.18893833333333
$ws.Range("N15").Value = 117.566815
$ws.Range("O15").Value = 0.637089375439711
$ws.Range("P15").Value = 0.637089375439711
$ws.Range("Q15").Value = 472.0274311652417
$ws.Range("R15").Value = 4248.246880487175
$ws.Range("S15").Value = 0.03928573676074666
$ws.Range("T15").Value = 0.03928573676074667
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 12.044915
$ws.Range("H16").Value = 36.134745
$ws.Range("I16").Value = 0.06166440420330686
$ws.Range("J16").Value = 0.06166440420330688
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.067472
$ws.Range("N16").Value = 18.202416
$ws.Range("O16").Value = 0.09863808797519777
$ws.Range("P16").Value = 0.09863808797519778
$ws.Range("Q16").Value = 73.08218450488
$ws.Range("R16").Value = 657.73966054392
$ws.Range("S16").Value = 0.006082458926743937
$ws.Range("T16").Value = 0.00608245892674394
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 12.044915
$ws.Range("H17").Value = 36.134745
$ws.Range("I17").Value = 0.06166440420330686
$ws.Range("J17").Value = 0.06166440420330688
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 10.968127
$ws.Range("N17").Value = 32.904381
$ws.Range("O17").Value = 0.1783073866594097
$ws.Range("P17").Value = 0.1783073866594098
$ws.Range("Q17").Value = 132.110157424205
$ws.Range("R17").Value = 1188.991416817845
$ws.Range("S17").Value = 0.01099521876340117
$ws.Range("T17").Value = 0.01099521876340117
